$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the existing row 497 (old rows 497..594
# shift down to 499..596), then fill the two new rows with their values.
$ws.Rows.Item(497).Resize(2).Insert()

# New row 497 (Primera / Crespo record / Repollo)
$ws.Range("A497").Value = 3
$ws.Range("B497").Value = "Femacal de La Calera"
$ws.Range("C497").Value = "Coquimbo"
$ws.Range("D497").Value = 44694
$ws.Range("E497").Value = 5
$ws.Range("F497").Value = 100112006
$ws.Range("G497").Value = "Repollo"
$ws.Range("H497").Value = "Crespo record"
$ws.Range("I497").Value = "Primera"
$ws.Range("J497").Value = 1800
$ws.Range("K497").Value = 1200
$ws.Range("L497").Value = 1300
$ws.Range("M497").Value = 1253
$ws.Range("N497").Value = "`$/unidad"
$ws.Range("O497").Value = "Provincia de Quillota"
$ws.Range("P497").Value = 1253
$ws.Range("Q497").Value = 1
$ws.Range("R497").Value = "Hortaliza"

# New row 498 (Segunda / Crespo record / Repollo)
$ws.Range("A498").Value = 3
$ws.Range("B498").Value = "Femacal de La Calera"
$ws.Range("C498").Value = "Coquimbo"
$ws.Range("D498").Value = 44694
$ws.Range("E498").Value = 5
$ws.Range("F498").Value = 100112006
$ws.Range("G498").Value = "Repollo"
$ws.Range("H498").Value = "Crespo record"
$ws.Range("I498").Value = "Segunda"
$ws.Range("J498").Value = 1800
$ws.Range("K498").Value = 900
$ws.Range("L498").Value = 900
$ws.Range("M498").Value = 900
$ws.Range("N498").Value = "`$/unidad"
$ws.Range("O498").Value = "Provincia de Quillota"
$ws.Range("P498").Value = 900
$ws.Range("Q498").Value = 1
$ws.Range("R498").Value = "Hortaliza"
